$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1867.8
$ws.Range("I32").Value = 1813.3334
$ws.Range("J32").Value = 1949.5
$ws.Range("K32").Value = 1813.3334
$ws.Range("L32").Value = 1949.5
$ws.Range("M32").Value = -1487.3334
$ws.Range("N32").Value = -2601.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7499.643
$ws.Range("I40").Value = 2999
$ws.Range("J40").Value = 9299.9
$ws.Range("K40").Value = 2999
$ws.Range("L40").Value = 9299.9
$ws.Range("M40").Value = -2824
$ws.Range("N40").Value = -9649.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4233.1665
$ws.Range("I74").Value = 4239.8
$ws.Range("J74").Value = 4200
$ws.Range("K74").Value = 4239.8
$ws.Range("L74").Value = 4200
$ws.Range("M74").Value = -3303.8
$ws.Range("N74").Value = -6072

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 4233.1665
$ws.Range("I77").Value = 4239.8
$ws.Range("J77").Value = 4200
$ws.Range("K77").Value = 21199
$ws.Range("L77").Value = 21000
$ws.Range("M77").Value = -16519
$ws.Range("N77").Value = -30360

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2726.8572
$ws.Range("I100").Value = 2876
$ws.Range("J100").Value = 1832
$ws.Range("K100").Value = 2876
$ws.Range("L100").Value = 1832
$ws.Range("M100").Value = -2335
$ws.Range("N100").Value = -2914

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 13364.613
$ws.Range("I138").Value = 12100
$ws.Range("J138").Value = 13406.767
$ws.Range("K138").Value = 36300
$ws.Range("L138").Value = 40220.301
$ws.Range("M138").Value = -31160
$ws.Range("N138").Value = -50500.301

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9117.473
$ws.Range("I32").Value = 8663.714
$ws.Range("J32").Value = 24999
$ws.Range("K32").Value = 8663.714
$ws.Range("L32").Value = 24999
$ws.Range("M32").Value = -8376.714
$ws.Range("N32").Value = -25573

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2330.3333
$ws.Range("I45").Value = 1595.5
$ws.Range("J45").Value = 3800
$ws.Range("K45").Value = 1595.5
$ws.Range("L45").Value = 3800
$ws.Range("M45").Value = -1218.5
$ws.Range("N45").Value = -4554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1708.2916
$ws.Range("I132").Value = 1319.3889
$ws.Range("J132").Value = 2875
$ws.Range("K132").Value = 3958.1667
$ws.Range("L132").Value = 8625
$ws.Range("M132").Value = -1428.1667
$ws.Range("N132").Value = -13685

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4529.4
$ws.Range("I105").Value = 3794.1
$ws.Range("J105").Value = 6000
$ws.Range("K105").Value = 3794.1
$ws.Range("L105").Value = 6000
$ws.Range("M105").Value = -2047.1
$ws.Range("N105").Value = -9494

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1312.1666
$ws.Range("I134").Value = 1024.6
$ws.Range("J134").Value = 2750
$ws.Range("K134").Value = 3073.8
$ws.Range("L134").Value = 8250
$ws.Range("M134").Value = -538.7999999999997
$ws.Range("N134").Value = -13320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 32561.334
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 32561.334
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 32561.334
$ws.Range("N88").Value = -33373.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 32561.334
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 32561.334
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 32561.334
$ws.Range("N91").Value = -35369.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5491.7856
$ws.Range("I99").Value = 3111.625
$ws.Range("J99").Value = 8665.333000000001
$ws.Range("K99").Value = 3111.625
$ws.Range("L99").Value = 8665.333000000001
$ws.Range("M99").Value = -1613.625
$ws.Range("N99").Value = -11661.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5491.7856
$ws.Range("I126").Value = 3111.625
$ws.Range("J126").Value = 8665.333000000001
$ws.Range("K126").Value = 9334.875
$ws.Range("L126").Value = 25995.999
$ws.Range("M126").Value = -6864.875
$ws.Range("N126").Value = -30935.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5222.1113
$ws.Range("I132").Value = 4542.7144
$ws.Range("J132").Value = 7600
$ws.Range("K132").Value = 13628.1432
$ws.Range("L132").Value = 22800
$ws.Range("M132").Value = -11098.1432
$ws.Range("N132").Value = -27860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3998.5
$ws.Range("I22").Value = 3998
$ws.Range("J22").Value = 3999
$ws.Range("K22").Value = 11994
$ws.Range("L22").Value = 11997
$ws.Range("M22").Value = -11825
$ws.Range("N22").Value = -12335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 41.8
$ws.Range("I26").Value = 41.8
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 125.4
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = 162.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 3998.5
$ws.Range("I27").Value = 3998
$ws.Range("J27").Value = 3999
$ws.Range("K27").Value = 11994
$ws.Range("L27").Value = 11997
$ws.Range("M27").Value = -11892
$ws.Range("N27").Value = -12201

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6548.846
$ws.Range("I122").Value = 6677.9165
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 20033.7495
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -17583.7495
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1262.9286
$ws.Range("I16").Value = 1257.7693
$ws.Range("J16").Value = 1330
$ws.Range("K16").Value = 1257.7693
$ws.Range("L16").Value = 1330
$ws.Range("M16").Value = -1087.7693
$ws.Range("N16").Value = -1670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1364
$ws.Range("I46").Value = 1296.5
$ws.Range("J46").Value = 1499
$ws.Range("K46").Value = 1296.5
$ws.Range("L46").Value = 1499
$ws.Range("M46").Value = -1108.5
$ws.Range("N46").Value = -1875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1053.6428
$ws.Range("I55").Value = 374.83334
$ws.Range("J55").Value = 1562.75
$ws.Range("K55").Value = 374.83334
$ws.Range("L55").Value = 1562.75
$ws.Range("M55").Value = -201.83334
$ws.Range("N55").Value = -1908.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 10000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 10000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10450

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 10000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 10000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3000
$ws.Range("N68").Value = -4498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 15000
$ws.Range("N71").Value = -22488

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4244.143
$ws.Range("I132").Value = 1941.8
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 5825.4
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -3295.4
$ws.Range("N132").Value = -35060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 12499
$ws.Range("I32").Value = 9998
$ws.Range("J32").Value = 15000
$ws.Range("K32").Value = 9998
$ws.Range("L32").Value = 15000
$ws.Range("M32").Value = -9681
$ws.Range("N32").Value = -15634

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6971545
$ws.Range("I100").Value = 9957722
$ws.Range("J100").Value = 3799.6667
$ws.Range("K100").Value = 19915444
$ws.Range("L100").Value = 7599.3334
$ws.Range("M100").Value = -19914903
$ws.Range("N100").Value = -8681.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2561.5
$ws.Range("I132").Value = 1597.8182
$ws.Range("J132").Value = 6095
$ws.Range("K132").Value = 4793.4546
$ws.Range("L132").Value = 18285
$ws.Range("M132").Value = -2263.4546
$ws.Range("N132").Value = -23345

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 34570.812
$ws.Range("I136").Value = 41010.31
$ws.Range("J136").Value = 6666.3335
$ws.Range("K136").Value = 123030.93
$ws.Range("L136").Value = 19999.0005
$ws.Range("M136").Value = -120480.93
$ws.Range("N136").Value = -25099.0005
